# Updated cryptos list on Fri Mar 22 14:24:50 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceCell($addr, $text) {
    # Price column values look numeric (e.g. "552.77", "1.00") and Excel's
    # COM layer auto-coerces such strings to floating point numbers, which
    # both changes the stored type and introduces binary rounding noise
    # (e.g. 552.76999999999998). Force the cell to Text format first so the
    # literal string is preserved exactly, then drop back to the default
    # "Normal" style so no stray number-format style sticks to the cell.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

function Set-TextCell($addr, $text) {
    $ws.Range($addr).Value = $text
}

# Row 2 - Bitcoin
Set-PriceCell "D2" "63.186.05"
Set-TextCell  "E2" "  -5.99%  "

# Row 3 - Ethereum
Set-PriceCell "D3" "3.330.82"
Set-TextCell  "E3" "  -6.24%  "

# Row 4 - TetherUSD
Set-PriceCell "D4" "1.00"
Set-TextCell  "E4" "  +0.13%  "

# Row 5 - BNB
Set-PriceCell "D5" "552.77"
Set-TextCell  "E5" "  -2.00%  "

# Row 6 - Solana
Set-PriceCell "D6" "169.12"
Set-TextCell  "E6" "  -10.86%  "

# Row 7 - XRP
Set-PriceCell "D7" "0.604"
Set-TextCell  "E7" "  -3.17%  "

# Row 8 - USDC
Set-PriceCell "D8" "1.00"
Set-TextCell  "E8" "  +0.23%  "

# Row 9 - Cardano
Set-PriceCell "D9" "0.603"
Set-TextCell  "E9" "  -5.22%  "

# Row 10 - Dogecoin
Set-PriceCell "D10" "0.147"
Set-TextCell  "E10" "  -4.10%  "

# Row 11 - Avalanche
Set-PriceCell "D11" "54.03"
Set-TextCell  "E11" "  -1.60%  "

# Row 12 - ShibaInu
Set-PriceCell "D12" "0.0000261"
Set-TextCell  "E12" "  -4.52%  "

# Row 13 - Polkadot
Set-PriceCell "D13" "8.78"
Set-TextCell  "E13" "  -6.37%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-PriceCell "D14" "3.870.34"
Set-TextCell  "E14" "  -5.60%  "

# Row 15 - TRON
Set-PriceCell "D15" "0.116"
Set-TextCell  "E15" "  -4.42%  "

# Row 16 - WrappedEther
Set-PriceCell "D16" "3.308.60"
Set-TextCell  "E16" "  -6.48%  "

# Row 17 - Chainlink
Set-PriceCell "D17" "17.50"
Set-TextCell  "E17" "  -6.08%  "

# Row 18 - WrappedBTC
Set-PriceCell "D18" "63.162.23"
Set-TextCell  "E18" "  -5.94%  "

# Row 19 - Uniswap
Set-PriceCell "D19" "11.43"
Set-TextCell  "E19" "  -5.77%  "

# Row 20 - Polygon
Set-PriceCell "D20" "0.960"
Set-TextCell  "E20" "  -4.14%  "

# Row 21 - BitcoinCash
Set-PriceCell "D21" "394.73"
Set-TextCell  "E21" "  -6.94%  "

# Row 22 - PancakeSwap
Set-PriceCell "D22" "4.01"
Set-TextCell  "E22" "  -2.93%  "

# Row 23 - Toncoin
Set-PriceCell "D23" "4.25"
Set-TextCell  "E23" "  +1.72%  "

# Row 24 - was Litecoin, now InternetComputer(DFINITY)
Set-TextCell  "B24" "InternetComputer(DFINITY)"
Set-TextCell  "C24" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-PriceCell "D24" "13.00"
Set-TextCell  "E24" "  +5.45%  "

# Row 25 - was InternetComputer(DFINITY), now Litecoin
Set-TextCell  "B25" "Litecoin"
Set-TextCell  "C25" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-PriceCell "D25" "81.05"
Set-TextCell  "E25" "  -5.12%  "

# Row 26 - RenderToken
Set-PriceCell "D26" "10.64"
Set-TextCell  "E26" "  -4.08%  "

# Row 27 - ImmutableX
Set-PriceCell "D27" "2.69"
Set-TextCell  "E27" "  -8.54%  "

# Row 28 - Filecoin
Set-PriceCell "D28" "8.51"
Set-TextCell  "E28" "  -6.51%  "

# Row 29 - EthereumClassic
Set-PriceCell "D29" "28.78"
Set-TextCell  "E29" "  -5.73%  "

# Row 30 - NEARProtocol
Set-PriceCell "D30" "6.40"
Set-TextCell  "E30" "  -4.08%  "

# Row 31 - Bittensor
Set-PriceCell "D31" "585.32"
Set-TextCell  "E31" "  -7.68%  "

# Row 32 - Cosmos
Set-PriceCell "D32" "11.17"
Set-TextCell  "E32" "  -5.17%  "

# Row 33 - Hedera
Set-PriceCell "D33" "0.104"
Set-TextCell  "E33" "  -6.82%  "

# Row 34 - OKB
Set-PriceCell "D34" "57.69"
Set-TextCell  "E34" "  -4.67%  "

# Row 35 - Kaspa (price unchanged)
Set-TextCell  "E35" "  -0.24%  "

# Row 36 - Dai (price unchanged)
Set-TextCell  "E36" "  +0.03%  "

# Row 37 - was InjectiveProtocol, now Stacks
Set-TextCell  "B37" "Stacks"
Set-TextCell  "C37" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-PriceCell "D37" "3.46"
Set-TextCell  "E37" "  +2.82%  "

# Row 38 - was Stacks, now InjectiveProtocol
Set-TextCell  "B38" "InjectiveProtocol"
Set-TextCell  "C38" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-PriceCell "D38" "35.47"
Set-TextCell  "E38" "  -7.86%  "

# Row 39 - Maker
Set-PriceCell "D39" "3.137.08"
Set-TextCell  "E39" "  -0.15%  "

# Row 40 - was TheGraph, now PEPE
Set-TextCell  "B40" "PEPE"
Set-TextCell  "C40" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-PriceCell "D40" "0.0₃0724"
Set-TextCell  "E40" "  -12.54%  "

# Row 41 - was PEPE, now TheGraph
Set-TextCell  "B41" "TheGraph"
Set-TextCell  "C41" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-PriceCell "D41" "0.361"
Set-TextCell  "E41" "  -7.20%  "

# Row 42 - FirstDigitalUSD
Set-PriceCell "D42" "0.997"
Set-TextCell  "E42" "  -0.35%  "

# Row 43 - ThetaToken
Set-PriceCell "D43" "2.77"
Set-TextCell  "E43" "  -3.85%  "

# Row 44 - Fetch.AI
Set-PriceCell "D44" "2.45"
Set-TextCell  "E44" "  -7.80%  "

# Row 45 - ApeXProtocol
Set-PriceCell "D45" "3.17"
Set-TextCell  "E45" "  -5.72%  "

# Row 46 - VeChain (price unchanged)
Set-TextCell  "E46" "  -5.14%  "

# Row 47 - WEMIXToken (price unchanged)
Set-TextCell  "E47" "  -6.20%  "

# Row 48 - Stellar
Set-PriceCell "D48" "0.126"
Set-TextCell  "E48" "  -4.82%  "

# Row 49 - Monero
Set-PriceCell "D49" "131.76"
Set-TextCell  "E49" "  -6.15%  "

# Row 50 - THORChain
Set-PriceCell "D50" "7.99"
Set-TextCell  "E50" "  -7.47%  "

# Row 51 - LidoDAOToken (price unchanged)
Set-TextCell  "E51" "  -1.09%  "
